$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 739.36365
$ws.Range("H112").Value = 1677.5526
$ws.Range("I112").Value = 919.6667
$ws.Range("J112").Value = 1819.6562
$ws.Range("K112").Value = 2759.0001
$ws.Range("L112").Value = 5458.9686
$ws.Range("M112").Value = -1651.0001
$ws.Range("N112").Value = -7674.9686
$ws.Range("H113").Value = 5445.364
$ws.Range("I113").Value = 2650
$ws.Range("J113").Value = 8799.799999999999
$ws.Range("K113").Value = 2650
$ws.Range("L113").Value = 8799.799999999999
$ws.Range("M113").Value = 604
$ws.Range("N113").Value = -15307.8
$ws.Range("H115").Value = 1543.25
$ws.Range("I115").Value = 803.0909
$ws.Range("J115").Value = 2169.5386
$ws.Range("K115").Value = 2409.2727
$ws.Range("L115").Value = 6508.6158
$ws.Range("M115").Value = -842.2727
$ws.Range("N115").Value = -9642.6158
$ws.Range("H118").Value = 834.4167
$ws.Range("I118").Value = 643.8
$ws.Range("J118").Value = 970.5714
$ws.Range("K118").Value = 1931.4
$ws.Range("L118").Value = 2911.7142
$ws.Range("M118").Value = -274.3999999999999
$ws.Range("N118").Value = -6225.7142
$ws.Range("H125").Value = 2429.5386
$ws.Range("I125").Value = 500
$ws.Range("J125").Value = 2590.3333
$ws.Range("K125").Value = 4500
$ws.Range("L125").Value = 23312.9997
$ws.Range("M125").Value = -2040
$ws.Range("N125").Value = -28232.9997
$ws.Range("H127").Value = 2006.5714
$ws.Range("I127").Value = 857
$ws.Range("J127").Value = 2256.4783
$ws.Range("K127").Value = 2571
$ws.Range("L127").Value = 6769.4349
$ws.Range("M127").Value = 2389
$ws.Range("N127").Value = -16689.4349
$ws.Range("H129").Value = 938.3111
$ws.Range("I129").Value = 495
$ws.Range("J129").Value = 958.93024
$ws.Range("K129").Value = 1485
$ws.Range("L129").Value = 2876.79072
$ws.Range("M129").Value = 3515
$ws.Range("N129").Value = -12876.79072
$ws.Range("H131").Value = 6590
$ws.Range("I131").Value = 12000
$ws.Range("J131").Value = 5237.5
$ws.Range("K131").Value = 36000
$ws.Range("L131").Value = 15712.5
$ws.Range("M131").Value = -30960
$ws.Range("N131").Value = -25792.5
$ws.Range("H132").Value = 378555.84
$ws.Range("I132").Value = 204863.95
$ws.Range("J132").Value = 2506281.5
$ws.Range("K132").Value = 614591.8500000001
$ws.Range("L132").Value = 7518844.5
$ws.Range("M132").Value = -612061.8500000001
$ws.Range("N132").Value = -7523904.5
$ws.Range("H135").Value = 274.73914
$ws.Range("I135").Value = 205.71428
$ws.Range("J135").Value = 999.5
$ws.Range("K135").Value = 1851.42852
$ws.Range("L135").Value = 8995.5
$ws.Range("M135").Value = 683.5714800000001
$ws.Range("N135").Value = -14065.5
$ws.Range("H137").Value = 2457.75
$ws.Range("I137").Value = 1216.9354
$ws.Range("J137").Value = 5416.615
$ws.Range("K137").Value = 3650.8062
$ws.Range("L137").Value = 16249.845
$ws.Range("M137").Value = -1100.8062
$ws.Range("N137").Value = -21349.845
$ws.Range("H138").Value = 1944.34
$ws.Range("I138").Value = 646.881
$ws.Range("J138").Value = 2883.8794
$ws.Range("K138").Value = 1940.643
$ws.Range("L138").Value = 8651.638199999999
$ws.Range("M138").Value = 3199.357
$ws.Range("N138").Value = -18931.6382
$ws.Range("H141").Value = 4968.346
$ws.Range("I141").Value = 5439.311
$ws.Range("J141").Value = 1940.7142
$ws.Range("K141").Value = 16317.933
$ws.Range("L141").Value = 5822.142599999999
$ws.Range("M141").Value = -11137.933
$ws.Range("N141").Value = -16182.1426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5726.9355
$ws.Range("I32").Value = 4491.2554
$ws.Range("J32").Value = 9598.733
$ws.Range("K32").Value = 4491.2554
$ws.Range("L32").Value = 9598.733
$ws.Range("M32").Value = -4204.2554
$ws.Range("H61").Value = 1117.3334
$ws.Range("I61").Value = 881.5625
$ws.Range("J61").Value = 3003.5
$ws.Range("K61").Value = 881.5625
$ws.Range("L61").Value = 3003.5
$ws.Range("M61").Value = -669.5625
$ws.Range("H63").Value = 7293199.5
$ws.Range("I63").Value = 9896064
$ws.Range("J63").Value = 5177.6
$ws.Range("K63").Value = 9896064
$ws.Range("L63").Value = 5177.6
$ws.Range("M63").Value = -9895378
$ws.Range("N63").Value = -6549.6
$ws.Range("H66").Value = 7293199.5
$ws.Range("I66").Value = 9896064
$ws.Range("J66").Value = 5177.6
$ws.Range("K66").Value = 49480320
$ws.Range("L66").Value = 25888
$ws.Range("M66").Value = -49476888
$ws.Range("N66").Value = -32752
$ws.Range("H74").Value = 3042.9211
$ws.Range("I74").Value = 3259.8276
$ws.Range("J74").Value = 2344
$ws.Range("K74").Value = 3259.8276
$ws.Range("L74").Value = 2344
$ws.Range("M74").Value = -2385.8276
$ws.Range("N74").Value = -4092
$ws.Range("H77").Value = 3042.9211
$ws.Range("I77").Value = 3259.8276
$ws.Range("J77").Value = 2344
$ws.Range("K77").Value = 16299.138
$ws.Range("L77").Value = 11720
$ws.Range("M77").Value = -11931.138
$ws.Range("N77").Value = -20456
$ws.Range("H102").Value = 1259.2593
$ws.Range("I102").Value = 1126.3158
$ws.Range("J102").Value = 1575
$ws.Range("K102").Value = 1126.3158
$ws.Range("L102").Value = 1575
$ws.Range("M102").Value = 495.6841999999999
$ws.Range("N102").Value = -4819
$ws.Range("H122").Value = 1951.5769
$ws.Range("I122").Value = 1252.5385
$ws.Range("J122").Value = 2650.6155
$ws.Range("K122").Value = 3757.6155
$ws.Range("L122").Value = 7951.8465
$ws.Range("M122").Value = -1307.6155
$ws.Range("N122").Value = -12851.8465
$ws.Range("H132").Value = 1919.6346
$ws.Range("I132").Value = 1080.641
$ws.Range("J132").Value = 4436.615
$ws.Range("K132").Value = 3241.923
$ws.Range("L132").Value = 13309.845
$ws.Range("M132").Value = -711.9230000000002
$ws.Range("N132").Value = -18369.845
$ws.Range("H136").Value = 1117.3334
$ws.Range("I136").Value = 881.5625
$ws.Range("J136").Value = 3003.5
$ws.Range("K136").Value = 2644.6875
$ws.Range("L136").Value = 9010.5
$ws.Range("M136").Value = -94.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1407.4902
$ws.Range("I134").Value = 1012.3
$ws.Range("J134").Value = 2844.5454
$ws.Range("K134").Value = 3036.9
$ws.Range("L134").Value = 8533.636200000001
$ws.Range("M134").Value = -501.8999999999996
$ws.Range("N134").Value = -13603.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1381.1852
$ws.Range("I132").Value = 745.0454999999999
$ws.Range("J132").Value = 4180.2
$ws.Range("K132").Value = 2235.1365
$ws.Range("L132").Value = 12540.6
$ws.Range("M132").Value = 294.8635000000004
$ws.Range("N132").Value = -17600.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1347.8636
$ws.Range("I5").Value = 449.1579
$ws.Range("J5").Value = 2030.88
$ws.Range("K5").Value = 1347.4737
$ws.Range("L5").Value = 6092.64
$ws.Range("M5").Value = -1235.4737
$ws.Range("N5").Value = -6316.64
$ws.Range("H113").Value = 614
$ws.Range("I113").Value = 571.5
$ws.Range("J113").Value = 670.6667
$ws.Range("K113").Value = 1714.5
$ws.Range("L113").Value = 2012.0001
$ws.Range("M113").Value = 455.5
$ws.Range("N113").Value = -6352.0001
$ws.Range("H117").Value = 4098.125
$ws.Range("I117").Value = 465.2
$ws.Range("J117").Value = 10153
$ws.Range("K117").Value = 1395.6
$ws.Range("L117").Value = 30459
$ws.Range("M117").Value = 2046.4
$ws.Range("N117").Value = -37343
$ws.Range("H122").Value = 3021.9119
$ws.Range("I122").Value = 673.2222
$ws.Range("J122").Value = 3867.44
$ws.Range("K122").Value = 6058.999800000001
$ws.Range("L122").Value = 34806.96
$ws.Range("M122").Value = -3608.999800000001
$ws.Range("N122").Value = -39706.96
$ws.Range("H131").Value = 9260061
$ws.Range("I131").Value = 55555870
$ws.Range("J131").Value = 898.55554
$ws.Range("K131").Value = 166667610
$ws.Range("L131").Value = 2695.66662
$ws.Range("M131").Value = -166662570
$ws.Range("N131").Value = -12775.66662
$ws.Range("H135").Value = 1347.8636
$ws.Range("I135").Value = 449.1579
$ws.Range("J135").Value = 2030.88
$ws.Range("K135").Value = 4042.4211
$ws.Range("L135").Value = 18277.92
$ws.Range("M135").Value = -1507.4211
$ws.Range("N135").Value = -23347.92
$ws.Range("H137").Value = 2589.2122
$ws.Range("I137").Value = 643.63635
$ws.Range("J137").Value = 3562
$ws.Range("K137").Value = 1930.90905
$ws.Range("L137").Value = 10686
$ws.Range("M137").Value = 3169.09095
$ws.Range("N137").Value = -20886

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2006.2693
$ws.Range("I102").Value = 1366
$ws.Range("J102").Value = 4695.4
$ws.Range("K102").Value = 1366
$ws.Range("L102").Value = 4695.4
$ws.Range("M102").Value = 256
$ws.Range("H107").Value = 7408035.5
$ws.Range("I107").Value = 485.45456
$ws.Range("J107").Value = 27778798
$ws.Range("K107").Value = 485.45456
$ws.Range("L107").Value = 27778798
$ws.Range("M107").Value = 1434.54544
$ws.Range("N107").Value = -27782638
$ws.Range("H126").Value = 2167.46
$ws.Range("I126").Value = 2173.9285
$ws.Range("J126").Value = 1850.5
$ws.Range("K126").Value = 6521.7855
$ws.Range("L126").Value = 5551.5
$ws.Range("M126").Value = -4051.7855
$ws.Range("N126").Value = -10491.5
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4066.611
$ws.Range("I122").Value = 1887.5
$ws.Range("J122").Value = 5809.9
$ws.Range("K122").Value = 5662.5
$ws.Range("L122").Value = 17429.7
$ws.Range("M122").Value = -3212.5
$ws.Range("N122").Value = -22329.7
$ws.Range("H125").Value = 41818.57
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 41818.57
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 41818.57
$ws.Range("N125").Value = -51658.57
$ws.Range("H136").Value = 2107.111
$ws.Range("I136").Value = 1178.5333
$ws.Range("J136").Value = 6750
$ws.Range("K136").Value = 3535.5999
$ws.Range("L136").Value = 20250
$ws.Range("M136").Value = -985.5999000000002
$ws.Range("N136").Value = -25350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 465.6875
$ws.Range("I113").Value = 390.2
$ws.Range("J113").Value = 591.5
$ws.Range("K113").Value = 1170.6
$ws.Range("L113").Value = 1774.5
$ws.Range("M113").Value = 999.4000000000001
$ws.Range("N113").Value = -6114.5
$ws.Range("H132").Value = 3789084.5
$ws.Range("I132").Value = 799.75
$ws.Range("J132").Value = 20836366
$ws.Range("K132").Value = 2399.25
$ws.Range("L132").Value = 62509098
$ws.Range("M132").Value = 130.75
$ws.Range("H136").Value = 2248.5103
$ws.Range("I136").Value = 550
$ws.Range("J136").Value = 6098.467
$ws.Range("K136").Value = 1650
$ws.Range("L136").Value = 18295.401
$ws.Range("M136").Value = 900
$ws.Range("N136").Value = -23395.401
$ws.Range("H138").Value = 42904.8
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 42904.8
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 42904.8
$ws.Range("N138").Value = -53184.8
